$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A93").Value = 1.598
$ws.Range("B93").Value = 1.091
$ws.Range("C93").Value = 2.342
